$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily expense rows appended below the existing data (rows 8-10)
$ws.Range("A8").Value = 43796
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 3

$ws.Range("A9").Value = 43797
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 2.5
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 25
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3

$ws.Range("A10").Value = 43798
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 2

# Match the date number format used in column A for the existing rows
$ws.Range("A8:A10").NumberFormat = "m/d/yy"

# Update selection to match the saved workbook state
$ws.Range("A11").Select() | Out-Null
